$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    -1,
    3.9999761565560537,
    4.3999737722116592,
    4.8399711494328255,
    5.3239682643761093,
    5.8563650908137213,
    6.442016267429997,
    7.0862178941729983,
    7.7948396835902987,
    8.574323651949328,
    9.4317560171442629,
    10.374946218188372,
    11.412440840007209,
    12.553684924007932,
    13.809053416408727,
    15.189958758049601,
    16.70895463385456,
    18.379866081825796
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 3).Value = $values[$i]
}

$ws.Columns.Item(3).ColumnWidth = 10.75
